# Apply the "contact info" formatting/text tweaks on the last slide
# (Slide 13 - "Content Placeholder 10") and the single rPr "dirty"
# touch-up in the Data Overview table on Slide 4.

$p = $ppt.ActivePresentation

# =====================================================================
# Slide 13 : Content Placeholder 10 (Email / LinkedIn / Github block)
# =====================================================================
$s13   = $p.Slides.Item(13)
$shape = $s13.Shapes.Item("Content Placeholder 10")
$tr    = $shape.TextFrame.TextRange

# --- Paragraph 1 : "Email: sivaprasad121333@gmail.com<spaces>ID No: S8064"
$para1 = $tr.Paragraphs(1)

# run 2 = the e-mail address -> bump to 18pt
$emailRun = $para1.Runs(2)
$emailRun.Font.Size = 18

# run 3 = the run of spaces that pads out to "ID No:" -> bump to 18pt and
# lengthen it (57 spaces -> 97 spaces) since the bigger font needs more
# padding to keep "ID No: S8064" in the same place.
$spacesRun = $para1.Runs(3)
$spacesRun.Font.Size = 18
$spacesRun.Text = "                                                                                                 "

# --- Paragraph 2 : "LinkedIn : https://www.linkedin.com/in/sivaprasad121333"
$para2 = $tr.Paragraphs(2)

# run 2 = the LinkedIn URL -> bump to 18pt
$linkedinUrlRun = $para2.Runs(2)
$linkedinUrlRun.Font.Size = 18

# --- Paragraph 3 : "Github: https://github.com/SIVAPRASAD121333"
$para3 = $tr.Paragraphs(3)

# run 2 currently holds ": https://github.com/SIVAPRASAD121333". Replace
# its text with ": " + the new (longer) repo URL, then carve the URL
# portion back out into its own run sized at 16pt (leaving ": " at the
# original/inherited size).
$githubRun = $para3.Runs(2)
$githubRun.Text = ": https://github.com/SIVAPRASAD121333/AMAZON-SALES-DATA-ANALYSIS"

$ghStart = $githubRun.Start
$ghLen   = $githubRun.Length

# skip the leading ": " (2 characters) and grab the rest of the run -
# setting a differing Font property on this sub-range splits it into a
# new run automatically.
$urlPart = $tr.Characters($ghStart + 2, $ghLen - 2)
$urlPart.Font.Size = 16

# =====================================================================
# Slide 4 : "Data Overview" table - row 13 ("time"), column 2
# (Description) - re-touch the run so it matches the saved-from-
# PowerPoint state.
# =====================================================================
$s4        = $p.Slides.Item(4)
$tblShape  = $s4.Shapes.Item("Content Placeholder 5")
$table     = $tblShape.Table
$descCell  = $table.Cell(13, 2)
$descRange = $descCell.Shape.TextFrame.TextRange
$descRange.Text = "The time at which the purchase was made"

Write-Host "Slide 13 contact block and Slide 4 table cell updated."
